$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PSSM score matrix for range B2:K21 (20 rows x 10 columns).
# Replaces the previous scores with the values from the supplemental figures update.
$data = @(
    @(-19.28024517505201, 1.998887576508555, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, 1.981021657283861, 1.668604674701444, -19.28024517505201, 3.391357437148838, -19.28024517505201, 1.796636298767654, -19.28024517505201, 1.048251165993111, -19.28024517505201),
    @(-19.28024517505201, 1.673471723397229, -19.28024517505201, -19.28024517505201, -19.28024517505201, 2.827723087318827, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(2.449548556875648, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, 1.832210240425431, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(3.861709297261871, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, 2.268481268625873),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, 2.912373030350917, -19.28024517505201, 2.948884903434795, -19.28024517505201, -19.28024517505201, -19.28024517505201, 2.039590644712689),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, 2.486390891937085, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, 1.609669437173073, 1.829343845313374),
    @(-19.28024517505201, -19.28024517505201, 1.53830345723229, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, 1.881982020694657),
    @(-19.28024517505201, -19.28024517505201, 1.722565034419781, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, 1.978608932042462, -19.28024517505201),
    @(-19.28024517505201, 2.105089061462546, 1.821325978561688, -19.28024517505201, -19.28024517505201, -19.28024517505201, 1.945129976175505, 4.321925942272439, 2.432834609051881, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, -19.28024517505201, 1.908636403469479, -19.28024517505201, 2.469749017266928, -19.28024517505201),
    @(-19.28024517505201, -19.28024517505201, 2.067417821428379, -19.28024517505201, -19.28024517505201, -19.28024517505201, 1.594682907181209, -19.28024517505201, -19.28024517505201, -19.28024517505201),
    @(-19.28024517505201, 1.043685883392121, 1.53143760496834, -19.28024517505201, 3.248983517383111, -19.28024517505201, 1.769431241297103, -19.28024517505201, -19.28024517505201, 1.937708382515175),
    @(-19.28024517505201, 1.322961885163876, -19.28024517505201, 1.72545954661237, -19.28024517505201, 2.372652739434164, 1.316808036410958, -19.28024517505201, -19.28024517505201, -19.28024517505201)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
